$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "Make '1st' conditional in external reference"
#
# Paragraph (table cell) that used to read:
#   Defendant ref: <<respondentExternalReference>><<cs_{!isBlank(respondent2ExternalReference)}>>
#
# becomes:
#   <<cs_{!isBlank(respondent2ExternalReference)}>>1st <<es_>>Defendant ref: <<respondentExternalReference>><<cs_{!isBlank(respondent2ExternalReference)}>>
#
# i.e. a new "<<cs_{!isBlank(respondent2ExternalReference)}>>1st <<es_>>"
# prefix is added before "Defendant ref: ..." (so "1st" only shows when
# there is a 2nd defendant), and the trailing conditional tag is kept
# (now uniformly coloured) and wrapped in a pair of OLE_LINK
# bookmarks, matching the ones already used around the equivalent
# "cs_{!isBlank(" tag used for the applicant/solicitor block earlier
# in the document.
# ------------------------------------------------------------------

$oldText = "Defendant ref: <<respondentExternalReference>><<cs_{!isBlank(respondent2ExternalReference)}>>"
$newText = "<<cs_{!isBlank(respondent2ExternalReference)}>>1st <<es_>>Defendant ref: <<respondentExternalReference>><<cs_{!isBlank(respondent2ExternalReference)}>>"

$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)

if (-not $found) {
    Write-Host "ERROR: target paragraph text not found"
}

$paraStart = $rng.Start
$paraEnd = $rng.Start + $newText.Length

# Make sure the whole replaced range is uniformly styled (the old,
# darker "0A0A0A" run colour used by the trailing conditional tag is
# gone in the new version - everything in this line is now "60686D").
$whole = $d.Range($paraStart, $paraEnd)
$whole.Font.Name = "GDSTransportWebsite"
$whole.Font.Color = 0x6D6860

# Wrap the final "<<cs_{!isBlank(respondent2ExternalReference)}>>" tag
# (the trailing one, right at the end of the paragraph) in the new
# OLE_LINK11 / OLE_LINK12 bookmark pair.
$tag = "<<cs_{!isBlank(respondent2ExternalReference)}>>"
$tagStart = $paraEnd - $tag.Length
$tagEnd = $paraEnd

$bmRange = $d.Range($tagStart, $tagEnd)
$d.Bookmarks.Add("OLE_LINK11", $bmRange)
$d.Bookmarks.Add("OLE_LINK12", $bmRange)
